# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data snapshot and the newly generated one.
# The same events are listed on both the "展览" sheet and the
# "全部类型" sheet (offset by one row there), so both need updating.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet (rows as in the original file)
$sheetExhibit.Range("F4").Value  = 174
$sheetExhibit.Range("F6").Value  = 5502
$sheetExhibit.Range("F9").Value  = 638
$sheetExhibit.Range("F11").Value = 1397
$sheetExhibit.Range("F12").Value = 34

# 全部类型 sheet (same events, shifted down by one row)
$sheetAll.Range("F4").Value  = 174
$sheetAll.Range("F7").Value  = 5502
$sheetAll.Range("F10").Value = 638
$sheetAll.Range("F12").Value = 1397
$sheetAll.Range("F13").Value = 34
